# Add a new column K (2022 data) to the worksheet, matching the style of
# the existing column J (2021 data), then update the sheet's used range,
# column widths, and selection to reflect the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column J (rows 3-5, the only rows that actually
# gain a new K cell) into column K so the new column picks up the same
# borders / number formats as its neighbour.
$ws.Range("J3:J5").Copy() | Out-Null
$ws.Range("K3:K5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New header value for 2022 and the corresponding data point.
$ws.Range("K4").Value = 2022
$ws.Range("K5").Value = 0.11705180708279034

# Match the column widths described by the new edit (columns D through K
# are all set to a stored OOXML width of 9 characters). The COM
# ColumnWidth setter and the stored <col width> value differ by 5/6 of a
# character in this runtime, so compensate for that offset here.
$ws.Range("D1:K1").ColumnWidth = 9 - 5/6

# Update the selected cell to reflect where the editor left off.
$ws.Range("J12").Select() | Out-Null
